# "Final Compile for Fall 2020 Class" - numerous tweaks made while
# creating help videos for students.
#
# Primary semantic change: the worksheet "Enrolment Statistics" is
# renamed to "Enrollment Statistics" (a spelling fix). Renaming the
# sheet is also what made this the active/selected tab when the file
# was last saved, and the embedded pie chart living on that sheet has
# its series formulas (which reference the sheet by name) updated to
# match the new name.

$wb = $excel.ActiveWorkbook

# 1) Rename the worksheet "Enrolment Statistics" -> "Enrollment Statistics".
$wsStats = $wb.Worksheets.Item("Enrolment Statistics")
$wsStats.Name = "Enrollment Statistics"

# 2) The embedded pie chart on that sheet (Chart 1) plots data from the
#    sheet via a SERIES() formula that hard-codes the sheet name, so it
#    needs to be repointed at the new sheet name as well.
$chartObj = $wsStats.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Enrollment Statistics'!`$B`$2,'Enrollment Statistics'!`$A`$3:`$A`$6,'Enrollment Statistics'!`$B`$3:`$B`$6,1)"

# 3) Make "Enrollment Statistics" the active/selected sheet (this is the
#    tab that was selected when the workbook was last saved).
$wsStats.Activate()
